$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 becomes a brand-new data point (most recent observation) ---
$ws.Range("D2").Value = 44860
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 14500
$ws.Range("P2").Value = 1115

# --- Rows 3..73 each inherit the values that used to sit one row above them ---
# Columns: Fecha(D, date serial), Volumen(J), Precio minimo(K), Precio maximo(L), Precio promedio ponderado(M), Precio $/Kg(P)
$rowsData = @(
    @(44644,400,15000,16000,15500,1192),
    @(44798,400,14000,15000,14500,1115),
    @(44792,400,14000,15000,14500,1115),
    @(44687,440,14000,15000,14500,1115),
    @(44309,400,26000,27000,26500,2038),
    @(44847,400,13500,14000,13750,1058),
    @(44333,440,24000,25000,24500,1885),
    @(44757,400,15000,16000,15500,1192),
    @(44839,520,14000,15000,14500,1115),
    @(44763,500,15000,16000,15500,1192),
    @(44846,400,13500,14000,13750,1058),
    @(44767,600,15000,16000,15500,1192),
    @(44764,400,15000,16000,15500,1192),
    @(44442,460,14000,15000,14500,1115),
    @(44708,440,13000,14000,13500,1038),
    @(44344,400,18500,19000,18750,1442),
    @(44410,600,14000,15000,14500,1115),
    @(44715,500,15000,16000,15500,1192),
    @(44740,500,16000,17000,16500,1269),
    @(44312,400,26000,27000,26500,2038),
    @(44761,480,14500,15000,14750,1135),
    @(44445,600,13000,14000,13500,1038),
    @(44418,500,14000,15000,14500,1115),
    @(44841,440,13500,14000,13750,1058),
    @(44400,600,15000,16000,15500,1192),
    @(44631,400,16000,17000,16500,1269),
    @(44599,400,15000,16000,15500,1192),
    @(44426,460,14000,15000,14500,1115),
    @(44658,400,15000,16000,15500,1192),
    @(44428,480,14000,15000,14500,1115),
    @(44754,480,15000,16000,15500,1192),
    @(44323,460,25000,26000,25500,1962),
    @(44837,600,14000,15000,14500,1115),
    @(44365,500,19500,20000,19750,1519),
    @(44704,400,13000,14000,13500,1038),
    @(44680,400,13500,14000,13750,1058),
    @(44750,480,15000,16000,15500,1192),
    @(44756,400,14500,15000,14750,1135),
    @(44414,500,14000,15000,14500,1115),
    @(44383,200,17000,18000,17500,1346),
    @(44692,400,14000,15000,14500,1115),
    @(44694,400,13000,14000,13500,1038),
    @(44831,600,14000,15000,14500,1115),
    @(44746,480,15000,16000,15500,1192),
    @(44701,440,14000,15000,14500,1115),
    @(44826,520,14000,15000,14500,1115),
    @(44419,600,14000,15000,14500,1115),
    @(44533,520,17000,18000,17500,1346),
    @(44505,400,16000,17000,16500,1269),
    @(44326,460,25000,26000,25500,1962),
    @(44657,460,15000,16000,15500,1192),
    @(44435,480,13000,14000,13500,1038),
    @(44412,600,14000,15000,14500,1115),
    @(44806,500,14000,15000,14500,1115),
    @(44771,480,14000,15000,14500,1115),
    @(44810,540,14000,15000,14500,1115),
    @(44582,520,15000,16000,15500,1192),
    @(44670,480,14500,15000,14750,1135),
    @(44747,440,15000,16000,15500,1192),
    @(44575,500,14000,15000,14500,1115),
    @(44855,500,13800,14000,13900,1069),
    @(44596,500,16000,17000,16500,1269),
    @(44803,520,14000,15000,14500,1115),
    @(44753,400,14500,15000,14750,1135),
    @(44498,400,14000,15000,14500,1115),
    @(44484,360,14000,15000,14500,1115),
    @(44799,460,14000,15000,14500,1115),
    @(44736,400,16000,17000,16500,1269),
    @(44379,600,17000,18000,17500,1346),
    @(44832,540,14000,15000,14500,1115),
    @(44335,480,24500,25000,24750,1904)
)

for ($i = 0; $i -lt $rowsData.Length; $i++) {
    $r = $i + 3
    $vals = $rowsData[$i]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("J$r").Value = $vals[1]
    $ws.Range("K$r").Value = $vals[2]
    $ws.Range("L$r").Value = $vals[3]
    $ws.Range("M$r").Value = $vals[4]
    $ws.Range("P$r").Value = $vals[5]
}

# Row 73 is a brand-new row (did not exist before): match the Fecha cells date number format to the rest of column D
$ws.Range("D73").NumberFormat = $ws.Range("D72").NumberFormat()

# Row 73 also needs the columns that are constant across every row in this sheet (same market / category / unit / origin / etc.)
$constCols = @("A","B","C","E","F","G","H","I","N","O","Q","R")
foreach ($col in $constCols) {
    $ws.Range("$col" + "73").Value = $ws.Range("$col" + "72").Value()
}
